$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input cells with new calibration values
$ws.Range("C9").Value = 3
$ws.Range("C10").Value = 2048
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 130

# Update the active selection on the sheet
$ws.Range("C15").Select()
